$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 29.223446
$ws.Cells.Item(2, 8).Value = 87.670338
$ws.Cells.Item(2, 9).Value = 0.0169041244192178
$ws.Cells.Item(2, 10).Value = 0.0169041244192178
$ws.Cells.Item(2, 13).Value = 3.135398666666667
$ws.Cells.Item(2, 14).Value = 9.406196000000001
$ws.Cells.Item(2, 15).Value = 0.1723049126704688
$ws.Cells.Item(2, 16).Value = 0.1723049126704688
$ws.Cells.Item(2, 17).Value = 91.62715362380534
$ws.Cells.Item(2, 18).Value = 824.6443826142481
$ws.Cells.Item(2, 19).Value = 0.002912663681824063
$ws.Cells.Item(2, 20).Value = 0.002912663681824063
$ws.Cells.Item(3, 7).Value = 29.223446
$ws.Cells.Item(3, 8).Value = 87.670338
$ws.Cells.Item(3, 9).Value = 0.0169041244192178
$ws.Cells.Item(3, 10).Value = 0.0169041244192178
$ws.Cells.Item(3, 15).Value = 0.1733096678828815
$ws.Cells.Item(3, 16).Value = 0.1733096678828815
$ws.Cells.Item(3, 17).Value = 92.16145562817199
$ws.Cells.Item(3, 18).Value = 829.4531006535479
$ws.Cells.Item(3, 19).Value = 0.002929648188945544
$ws.Cells.Item(3, 20).Value = 0.002929648188945544
$ws.Cells.Item(4, 7).Value = 29.223446
$ws.Cells.Item(4, 8).Value = 87.670338
$ws.Cells.Item(4, 9).Value = 0.0169041244192178
$ws.Cells.Item(4, 10).Value = 0.0169041244192178
$ws.Cells.Item(4, 13).Value = 0.4900660000000001
$ws.Cells.Item(4, 14).Value = 1.470198
$ws.Cells.Item(4, 15).Value = 0.02693143306797965
$ws.Cells.Item(4, 16).Value = 0.02693143306797965
$ws.Cells.Item(4, 17).Value = 14.321417287436
$ws.Cells.Item(4, 18).Value = 128.892755586924
$ws.Cells.Item(4, 19).Value = 0.0004552522953689646
$ws.Cells.Item(4, 20).Value = 0.0004552522953689646
$ws.Cells.Item(5, 7).Value = 29.223446
$ws.Cells.Item(5, 8).Value = 87.670338
$ws.Cells.Item(5, 9).Value = 0.0169041244192178
$ws.Cells.Item(5, 10).Value = 0.0169041244192178
$ws.Cells.Item(5, 13).Value = 11.417657
$ws.Cells.Item(5, 14).Value = 34.252971
$ws.Cells.Item(5, 15).Value = 0.62745398637867
$ws.Cells.Item(5, 16).Value = 0.6274539863786701
$ws.Cells.Item(5, 17).Value = 333.663282786022
$ws.Cells.Item(5, 18).Value = 3002.969545074198
$ws.Cells.Item(5, 19).Value = 0.01060656025307923
$ws.Cells.Item(5, 20).Value = 0.01060656025307923
$ws.Cells.Item(6, 9).Value = 0.9471112884046843
$ws.Cells.Item(6, 10).Value = 0.9471112884046842
$ws.Cells.Item(6, 13).Value = 3.135398666666667
$ws.Cells.Item(6, 14).Value = 9.406196000000001
$ws.Cells.Item(6, 15).Value = 0.1723049126704688
$ws.Cells.Item(6, 16).Value = 0.1723049126704688
$ws.Cells.Item(6, 17).Value = 5133.724135562876
$ws.Cells.Item(6, 18).Value = 46203.51722006589
$ws.Cells.Item(6, 19).Value = 0.1631919278377843
$ws.Cells.Item(6, 20).Value = 0.1631919278377843
$ws.Cells.Item(7, 9).Value = 0.9471112884046843
$ws.Cells.Item(7, 10).Value = 0.9471112884046842
$ws.Cells.Item(7, 15).Value = 0.1733096678828815
$ws.Cells.Item(7, 16).Value = 0.1733096678828815
$ws.Cells.Item(7, 19).Value = 0.1641435428415438
$ws.Cells.Item(7, 20).Value = 0.1641435428415438
$ws.Cells.Item(8, 9).Value = 0.9471112884046843
$ws.Cells.Item(8, 10).Value = 0.9471112884046842
$ws.Cells.Item(8, 13).Value = 0.4900660000000001
$ws.Cells.Item(8, 14).Value = 1.470198
$ws.Cells.Item(8, 15).Value = 0.02693143306797965
$ws.Cells.Item(8, 16).Value = 0.02693143306797965
$ws.Cells.Item(8, 17).Value = 802.4063028939934
$ws.Cells.Item(8, 18).Value = 7221.656726045941
$ws.Cells.Item(8, 19).Value = 0.02550706427159873
$ws.Cells.Item(8, 20).Value = 0.02550706427159872
$ws.Cells.Item(9, 9).Value = 0.9471112884046843
$ws.Cells.Item(9, 10).Value = 0.9471112884046842
$ws.Cells.Item(9, 13).Value = 11.417657
$ws.Cells.Item(9, 14).Value = 34.252971
$ws.Cells.Item(9, 15).Value = 0.62745398637867
$ws.Cells.Item(9, 16).Value = 0.6274539863786701
$ws.Cells.Item(9, 17).Value = 18694.62468541324
$ws.Cells.Item(9, 18).Value = 168251.6221687191
$ws.Cells.Item(9, 19).Value = 0.5942687534537574
$ws.Cells.Item(9, 20).Value = 0.5942687534537574
$ws.Cells.Item(10, 7).Value = 37.39212666666667
$ws.Cells.Item(10, 8).Value = 112.17638
$ws.Cells.Item(10, 9).Value = 0.02162924801792661
$ws.Cells.Item(10, 10).Value = 0.0216292480179266
$ws.Cells.Item(10, 13).Value = 3.135398666666667
$ws.Cells.Item(10, 14).Value = 9.406196000000001
$ws.Cells.Item(10, 15).Value = 0.1723049126704688
$ws.Cells.Item(10, 16).Value = 0.1723049126704688
$ws.Cells.Item(10, 17).Value = 117.2392240944978
$ws.Cells.Item(10, 18).Value = 1055.15301685048
$ws.Cells.Item(10, 19).Value = 0.003726825690856755
$ws.Cells.Item(10, 20).Value = 0.003726825690856755
$ws.Cells.Item(11, 7).Value = 37.39212666666667
$ws.Cells.Item(11, 8).Value = 112.17638
$ws.Cells.Item(11, 9).Value = 0.02162924801792661
$ws.Cells.Item(11, 10).Value = 0.0216292480179266
$ws.Cells.Item(11, 15).Value = 0.1733096678828815
$ws.Cells.Item(11, 16).Value = 0.1733096678828815
$ws.Cells.Item(11, 17).Value = 117.9228768103867
$ws.Cells.Item(11, 18).Value = 1061.30589129348
$ws.Cells.Item(11, 19).Value = 0.003748557790543333
$ws.Cells.Item(11, 20).Value = 0.003748557790543332
$ws.Cells.Item(12, 7).Value = 37.39212666666667
$ws.Cells.Item(12, 8).Value = 112.17638
$ws.Cells.Item(12, 9).Value = 0.02162924801792661
$ws.Cells.Item(12, 10).Value = 0.0216292480179266
$ws.Cells.Item(12, 13).Value = 0.4900660000000001
$ws.Cells.Item(12, 14).Value = 1.470198
$ws.Cells.Item(12, 15).Value = 0.02693143306797965
$ws.Cells.Item(12, 16).Value = 0.02693143306797965
$ws.Cells.Item(12, 17).Value = 18.32460994702667
$ws.Cells.Item(12, 18).Value = 164.92148952324
$ws.Cells.Item(12, 19).Value = 0.000582506645305522
$ws.Cells.Item(12, 20).Value = 0.0005825066453055219
$ws.Cells.Item(13, 7).Value = 37.39212666666667
$ws.Cells.Item(13, 8).Value = 112.17638
$ws.Cells.Item(13, 9).Value = 0.02162924801792661
$ws.Cells.Item(13, 10).Value = 0.0216292480179266
$ws.Cells.Item(13, 13).Value = 11.417657
$ws.Cells.Item(13, 14).Value = 34.252971
$ws.Cells.Item(13, 15).Value = 0.62745398637867
$ws.Cells.Item(13, 16).Value = 0.6274539863786701
$ws.Cells.Item(13, 17).Value = 426.9304767805534
$ws.Cells.Item(13, 18).Value = 3842.374291024981
$ws.Cells.Item(13, 19).Value = 0.013571357891221
$ws.Cells.Item(13, 20).Value = 0.013571357891221
$ws.Cells.Item(14, 7).Value = 24.817167
$ws.Cells.Item(14, 8).Value = 74.45150100000001
$ws.Cells.Item(14, 9).Value = 0.01435533915817136
$ws.Cells.Item(14, 10).Value = 0.01435533915817136
$ws.Cells.Item(14, 13).Value = 3.135398666666667
$ws.Cells.Item(14, 14).Value = 9.406196000000001
$ws.Cells.Item(14, 15).Value = 0.1723049126704688
$ws.Cells.Item(14, 16).Value = 0.1723049126704688
$ws.Cells.Item(14, 17).Value = 77.81171232224402
$ws.Cells.Item(14, 18).Value = 700.3054109001962
$ws.Cells.Item(14, 19).Value = 0.002473495460003678
$ws.Cells.Item(14, 20).Value = 0.002473495460003678
$ws.Cells.Item(15, 7).Value = 24.817167
$ws.Cells.Item(15, 8).Value = 74.45150100000001
$ws.Cells.Item(15, 9).Value = 0.01435533915817136
$ws.Cells.Item(15, 10).Value = 0.01435533915817136
$ws.Cells.Item(15, 15).Value = 0.1733096678828815
$ws.Cells.Item(15, 16).Value = 0.1733096678828815
$ws.Cells.Item(15, 17).Value = 78.26545285889399
$ws.Cells.Item(15, 18).Value = 704.389075730046
$ws.Cells.Item(15, 19).Value = 0.002487919061848802
$ws.Cells.Item(15, 20).Value = 0.002487919061848802
$ws.Cells.Item(16, 7).Value = 24.817167
$ws.Cells.Item(16, 8).Value = 74.45150100000001
$ws.Cells.Item(16, 9).Value = 0.01435533915817136
$ws.Cells.Item(16, 10).Value = 0.01435533915817136
$ws.Cells.Item(16, 13).Value = 0.4900660000000001
$ws.Cells.Item(16, 14).Value = 1.470198
$ws.Cells.Item(16, 15).Value = 0.02693143306797965
$ws.Cells.Item(16, 16).Value = 0.02693143306797965
$ws.Cells.Item(16, 17).Value = 12.162049763022
$ws.Cells.Item(16, 18).Value = 109.458447867198
$ws.Cells.Item(16, 19).Value = 0.0003866098557064394
$ws.Cells.Item(16, 20).Value = 0.0003866098557064393
$ws.Cells.Item(17, 7).Value = 24.817167
$ws.Cells.Item(17, 8).Value = 74.45150100000001
$ws.Cells.Item(17, 9).Value = 0.01435533915817136
$ws.Cells.Item(17, 10).Value = 0.01435533915817136
$ws.Cells.Item(17, 13).Value = 11.417657
$ws.Cells.Item(17, 14).Value = 34.252971
$ws.Cells.Item(17, 15).Value = 0.62745398637867
$ws.Cells.Item(17, 16).Value = 0.6274539863786701
$ws.Cells.Item(17, 17).Value = 283.353900517719
$ws.Cells.Item(17, 18).Value = 2550.185104659472
$ws.Cells.Item(17, 19).Value = 0.009007314780612442
$ws.Cells.Item(17, 20).Value = 0.009007314780612442
